# Auto-generated edit script: update Sheets via scheduled runner
# Applies per-cell numeric updates across ALC, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 6693.737
$ws.Range("I98").Value = 7761.3125
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 7761.3125
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = -6263.3125
$ws.Range("N98").Value = -3996
$ws.Range("H122").Value = 6693.737
$ws.Range("I122").Value = 7761.3125
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 23283.9375
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -20833.9375
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3148370.5
$ws.Range("I31").Value = 2425.5334
$ws.Range("J31").Value = 7251777
$ws.Range("K31").Value = 2425.5334
$ws.Range("L31").Value = 7251777
$ws.Range("M31").Value = -2130.5334
$ws.Range("N31").Value = -7252367
$ws.Range("H34").Value = 3148370.5
$ws.Range("I34").Value = 2425.5334
$ws.Range("J34").Value = 7251777
$ws.Range("K34").Value = 2425.5334
$ws.Range("L34").Value = 7251777
$ws.Range("M34").Value = -2223.5334
$ws.Range("N34").Value = -7252181
$ws.Range("H59").Value = 30705
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 30705
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 30705
$ws.Range("M59").Value = $null
$ws.Range("N59").Value = -32995
$ws.Range("H99").Value = 2304.4
$ws.Range("I99").Value = 2099.4546
$ws.Range("J99").Value = 2554.889
$ws.Range("K99").Value = 2099.4546
$ws.Range("L99").Value = 2554.889
$ws.Range("M99").Value = -601.4546
$ws.Range("N99").Value = -5550.889
$ws.Range("H126").Value = 2304.4
$ws.Range("I126").Value = 2099.4546
$ws.Range("J126").Value = 2554.889
$ws.Range("K126").Value = 6298.3638
$ws.Range("L126").Value = 7664.667
$ws.Range("M126").Value = -3828.3638
$ws.Range("N126").Value = -12604.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1180
$ws.Range("J34").Value = 1900
$ws.Range("L34").Value = 5700
$ws.Range("N34").Value = -5868
$ws.Range("H39").Value = 2000
$ws.Range("H46").Value = 200
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null
$ws.Range("H55").Value = 2654.5454
$ws.Range("J55").Value = 2654.5454
$ws.Range("L55").Value = 7963.6362
$ws.Range("N55").Value = -8317.636200000001
$ws.Range("H62").Value = 1987.5
$ws.Range("I62").Value = 1316.6666
$ws.Range("K62").Value = 3949.9998
$ws.Range("M62").Value = -3263.9998
$ws.Range("H65").Value = 1987.5
$ws.Range("I65").Value = 1316.6666
$ws.Range("K65").Value = 11849.9994
$ws.Range("M65").Value = -8417.999400000001
$ws.Range("H69").Value = 1663.6364
$ws.Range("I69").Value = 766.6667
$ws.Range("J69").Value = 2000
$ws.Range("K69").Value = 2300.0001
$ws.Range("L69").Value = 6000
$ws.Range("M69").Value = -1489.0001
$ws.Range("N69").Value = -7622
$ws.Range("H72").Value = 1663.6364
$ws.Range("I72").Value = 766.6667
$ws.Range("J72").Value = 2000
$ws.Range("K72").Value = 6900.0003
$ws.Range("L72").Value = 18000
$ws.Range("M72").Value = -2844.0003
$ws.Range("N72").Value = -26112
$ws.Range("H86").Value = 2550.5
$ws.Range("I86").Value = 500
$ws.Range("J86").Value = 2960.6
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 8881.799999999999
$ws.Range("M86").Value = -314
$ws.Range("N86").Value = -11253.8
$ws.Range("H89").Value = 2550.5
$ws.Range("I89").Value = 500
$ws.Range("J89").Value = 2960.6
$ws.Range("K89").Value = 4500
$ws.Range("L89").Value = 26645.4
$ws.Range("M89").Value = 1428
$ws.Range("N89").Value = -38501.39999999999
$ws.Range("H127").Value = 1314
$ws.Range("J127").Value = 1314
$ws.Range("L127").Value = 3942
$ws.Range("N127").Value = -13862
$ws.Range("H131").Value = 762.65
$ws.Range("I131").Value = 375
$ws.Range("J131").Value = 796.3587
$ws.Range("K131").Value = 1125
$ws.Range("L131").Value = 2389.0761
$ws.Range("M131").Value = 3915
$ws.Range("N131").Value = -12469.0761

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 47622412
$ws.Range("I122").Value = 166671300
$ws.Range("J122").Value = 2858.9333
$ws.Range("K122").Value = 500013900
$ws.Range("L122").Value = 8576.7999
$ws.Range("M122").Value = -500011450
$ws.Range("N122").Value = -13476.7999
$ws.Range("H126").Value = 2003.7222
$ws.Range("I126").Value = 3139.7144
$ws.Range("K126").Value = 9419.143199999999
$ws.Range("M126").Value = -6949.143199999999
$ws.Range("H127").Value = 28333.334
$ws.Range("J127").Value = 28333.334
$ws.Range("L127").Value = 28333.334
$ws.Range("N127").Value = -38253.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 100001420
$ws.Range("I7").Value = 1775
$ws.Range("J7").Value = 500000000
$ws.Range("K7").Value = 1775
$ws.Range("L7").Value = 500000000
$ws.Range("M7").Value = -1663
$ws.Range("N7").Value = -500000224
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = $null
$ws.Range("H100").Value = 41668350
$ws.Range("I100").Value = 66667930
$ws.Range("K100").Value = 66667930
$ws.Range("M100").Value = -66667389
$ws.Range("H122").Value = 2511.1428
$ws.Range("I122").Value = 1866
$ws.Range("J122").Value = 2995
$ws.Range("K122").Value = 5598
$ws.Range("L122").Value = 8985
$ws.Range("M122").Value = -3148
$ws.Range("N122").Value = -13885
$ws.Range("H126").Value = 100001420
$ws.Range("I126").Value = 1775
$ws.Range("J126").Value = 500000000
$ws.Range("K126").Value = 5325
$ws.Range("L126").Value = 1500000000
$ws.Range("M126").Value = -2855
$ws.Range("N126").Value = -1500004940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 50007344
$ws.Range("I126").Value = 62508496
$ws.Range("J126").Value = 2725
$ws.Range("K126").Value = 187525488
$ws.Range("L126").Value = 8175
$ws.Range("M126").Value = -187523018
$ws.Range("N126").Value = -13115

